# add dvd insert function
# Convert the RATING column (G) from Korean text labels to plain numeric
# age values, matching the new "insert" workflow that stores ratings as
# numbers instead of shared-string labels.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 15
$ws.Range("G4").Value = 19
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 15
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 15
$ws.Range("G9").Value = 12
$ws.Range("G10").Value = 12
$ws.Range("G11").Value = 12

# Leave the active selection where the user last clicked after entering data.
$ws.Range("K11").Select() | Out-Null
